$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.634.36'
$ws.Range("E2").Value = '  -1.49%  '
$ws.Range("D3").Value = '2.357.05'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.78%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -1.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.77'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0922'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.41'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.83%  '
$ws.Range("E13").Value = '  -4.68%  '
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("D16").Value = '2.708.72'
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").Value = '2.350.67'
$ws.Range("E17").Value = '  -3.50%  '
$ws.Range("E18").Value = '  +12.14%  '
$ws.Range("D19").Value = '42.470.71'
$ws.Range("E19").Value = '  -1.85%  '
$ws.Range("E20").Value = '  -2.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '76.17'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.70'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '268.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.55%  '
$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +11.83%  '
$ws.Range("B25").Value = 'ImmutableX'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.31'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -10.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.47'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.06%  '
$ws.Range("E29").Value = '  -2.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.59'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.10%  '
$ws.Range("E31").Value = '  -2.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0899'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.31'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -9.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.04'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.133'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.59'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -8.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0358'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.95'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.82%  '
$ws.Range("E39").Value = '  +1.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.80'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.46%  '
$ws.Range("E41").Value = '  +1.33%  '
$ws.Range("E42").Value = '  +0.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.15'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.32%  '
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '119.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.91%  '
$ws.Range("B46").Value = 'BitcoinSV'
$ws.Range("C46").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +32.21%  '
$ws.Range("E47").Value = '  -7.85%  '
$ws.Range("E48").Value = '  -2.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.34%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '1.564.23'
$ws.Range("E50").Value = '  +4.55%  '
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.26'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.59%  '
